# Commit: "Various fixes related to GC data and sample collection"
#
# The site U-01 (eval.status "PI" = Physically Inaccessible) on the "data"
# sheet never had any real field measurements recorded for it - remove
# that row entirely so the remaining rows (U-02, U-03, ...) shift up.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Rows.Item(3).Delete()

# Restore the user's on-screen selection/scroll position after the edit.
$ws.Range("N10").Select() | Out-Null
